# "Save before changing branches"
# Update the "params" sheet data table: Visibility, Human/Autonomous Acft
# counts, and Separation Distance columns, then refresh the UI selection
# on both sheets.

$wb = $excel.ActiveWorkbook
$runs = $wb.Worksheets.Item("runs")
$params = $wb.Worksheets.Item("params")

# --- params sheet cell edits -------------------------------------------------

# Column D (Visibility)
$params.Range("D2").Value = 0.8
$params.Range("D6").Value = 0.6
$params.Range("D11").Value = 0.6
$params.Range("D16").Value = 0.6

# Column B (Human Acft) / Column C (Autonomous Acft)
$params.Range("B7").Value = 10
$params.Range("B8").Value = 20
$params.Range("C8").Value = 10
$params.Range("B9").Value = 30
$params.Range("C9").Value = 10
$params.Range("B10").Value = 40
$params.Range("C10").Value = 10
$params.Range("B11").Value = 50
$params.Range("C11").Value = 10
$params.Range("C12").Value = 10
$params.Range("B13").Value = 20
$params.Range("B14").Value = 30
$params.Range("C14").Value = 10
$params.Range("B15").Value = 40
$params.Range("C15").Value = 10
$params.Range("B16").Value = 50
$params.Range("C16").Value = 10

# Column J (Separation Distance (m)) - whole block changes 250 -> 100
$params.Range("J2:J16").Value = 100

# --- selection / view state --------------------------------------------------

$runs.Activate()
$runs.Range("G16").Select()

$params.Activate()
$params.Range("I24").Select()
